$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix classification report sort order: swap the data rows for
# A <-> AAA, B <-> BBB, C <-> CCC so ratings are sorted correctly.

# Row 2 (was A) becomes AAA
$ws.Range("A2").Value = "AAA"
$ws.Range("B2").Value = "0.8261"
$ws.Range("C2").Value = "0.7917"
$ws.Range("D2").Value = "0.8085"
$ws.Range("E2").Value = "24"

# Row 4 (was AAA) becomes A
$ws.Range("A4").Value = "A"
$ws.Range("B4").Value = "0.5874"
$ws.Range("C4").Value = "0.6298"
$ws.Range("D4").Value = "0.6079"
$ws.Range("E4").Value = "208"

# Row 5 (was B) becomes BBB
$ws.Range("A5").Value = "BBB"
$ws.Range("B5").Value = "0.6909"
$ws.Range("C5").Value = "0.6281"
$ws.Range("D5").Value = "0.6580"
$ws.Range("E5").Value = "363"

# Row 7 (was BBB) becomes B
$ws.Range("A7").Value = "B"
$ws.Range("B7").Value = "0.5864"
$ws.Range("C7").Value = "0.7273"
$ws.Range("D7").Value = "0.6493"
$ws.Range("E7").Value = "154"

# Row 8 (was C) becomes CCC
$ws.Range("A8").Value = "CCC"
$ws.Range("B8").Value = "0.5556"
$ws.Range("C8").Value = "0.7692"
$ws.Range("D8").Value = "0.6452"
$ws.Range("E8").Value = "26"

# Row 10 (was CCC) becomes C
$ws.Range("A10").Value = "C"
$ws.Range("B10").Value = "1.0000"
$ws.Range("C10").Value = "1.0000"
$ws.Range("D10").Value = "1.0000"
$ws.Range("E10").Value = "4"
